# Fruta / hortaliza, semanal
#
# A new weekly price observation for Maracuya (Agricola del Norte S.A. de
# Arica) is inserted at the top of the series (row 44). Every existing
# observation in rows 44-72 shifts down by one row (-> 45-73), so the oldest
# observation (previously row 72) becomes the new row 73.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: shift existing rows 44-72 down into 45-73 (columns D,L,M,N,O,P,S) ---
# Listed bottom (row 72->73) to top (row 44->45) so that, as we write into a
# destination row, the still-to-be-read source rows above it are untouched.
$shiftRows = @(
    @{ Src = 72; D = 44284; L = "Primera"; M = 120; N = 20000; O = 21000; P = 20500; S = 1025 }
    @{ Src = 71; D = 44270; L = "Primera"; M = 120; N = 24000; O = 25000; P = 24500; S = 1225 }
    @{ Src = 70; D = 44270; L = "Especial"; M = 120; N = 26000; O = 27000; P = 26500; S = 1325 }
    @{ Src = 69; D = 44232; L = "Primera"; M = 120; N = 30000; O = 31000; P = 30500; S = 1525 }
    @{ Src = 68; D = 44232; L = "Especial"; M = 120; N = 32000; O = 33000; P = 32500; S = 1625 }
    @{ Src = 67; D = 44281; L = "Primera"; M = 120; N = 23000; O = 24000; P = 23500; S = 1175 }
    @{ Src = 66; D = 44281; L = "Especial"; M = 120; N = 24000; O = 25000; P = 24500; S = 1225 }
    @{ Src = 65; D = 44316; L = "Segunda"; M = 120; N = 17000; O = 18000; P = 17500; S = 875 }
    @{ Src = 64; D = 44316; L = "Primera"; M = 140; N = 19000; O = 20000; P = 19500; S = 975 }
    @{ Src = 63; D = 44344; L = "Segunda"; M = 120; N = 24000; O = 25000; P = 24500; S = 1225 }
    @{ Src = 62; D = 44344; L = "Primera"; M = 120; N = 26000; O = 27000; P = 26500; S = 1325 }
    @{ Src = 61; D = 44344; L = "Especial"; M = 120; N = 29000; O = 30000; P = 29500; S = 1475 }
    @{ Src = 60; D = 44309; L = "Segunda"; M = 120; N = 15000; O = 16000; P = 15500; S = 775 }
    @{ Src = 59; D = 44309; L = "Primera"; M = 120; N = 17000; O = 18000; P = 17500; S = 875 }
    @{ Src = 58; D = 44277; L = "Primera"; M = 120; N = 23000; O = 24000; P = 23500; S = 1175 }
    @{ Src = 57; D = 44277; L = "Especial"; M = 120; N = 24000; O = 25000; P = 24500; S = 1225 }
    @{ Src = 56; D = 44333; L = "Primera"; M = 120; N = 27000; O = 28000; P = 27500; S = 1375 }
    @{ Src = 55; D = 44333; L = "Especial"; M = 120; N = 29000; O = 30000; P = 29500; S = 1475 }
    @{ Src = 54; D = 44265; L = "Primera"; M = 140; N = 23000; O = 24000; P = 23500; S = 1175 }
    @{ Src = 53; D = 44292; L = "Primera"; M = 160; N = 17000; O = 18000; P = 17500; S = 875 }
    @{ Src = 52; D = 44389; L = "Primera"; M = 120; N = 24000; O = 25000; P = 24500; S = 1225 }
    @{ Src = 51; D = 44358; L = "Segunda"; M = 160; N = 22000; O = 23000; P = 22500; S = 1125 }
    @{ Src = 50; D = 44358; L = "Primera"; M = 120; N = 25000; O = 26000; P = 25500; S = 1275 }
    @{ Src = 49; D = 44351; L = "Especial"; M = 160; N = 29000; O = 30000; P = 29500; S = 1475 }
    @{ Src = 48; D = 44330; L = "Segunda"; M = 120; N = 21000; O = 22000; P = 21500; S = 1075 }
    @{ Src = 47; D = 44330; L = "Primera"; M = 120; N = 24000; O = 25000; P = 24500; S = 1225 }
    @{ Src = 46; D = 44330; L = "Especial"; M = 100; N = 27000; O = 28000; P = 27500; S = 1375 }
    @{ Src = 45; D = 44231; L = "Primera"; M = 100; N = 34000; O = 35000; P = 34500; S = 1725 }
    @{ Src = 44; D = 44260; L = "Primera"; M = 120; N = 21000; O = 22000; P = 21500; S = 1075 }
)

foreach ($row in $shiftRows) {
    $dest = $row.Src + 1
    $ws.Cells.Item($dest, 4).Value  = $row.D   # D: Fecha
    $ws.Cells.Item($dest, 12).Value = $row.L   # L: Categoría/Calidad
    $ws.Cells.Item($dest, 13).Value = $row.M   # M: Calibre/tamaño
    $ws.Cells.Item($dest, 14).Value = $row.N   # N: Precio mínimo
    $ws.Cells.Item($dest, 15).Value = $row.O   # O: Precio máximo
    $ws.Cells.Item($dest, 16).Value = $row.P   # P: Precio promedio
    $ws.Cells.Item($dest, 19).Value = $row.S   # S: Precio por kilo
}

# --- Step 2: row 73 is a brand-new row; copy the constant columns from row 72 ---
# (column D keeps the date-styled NumberFormat used throughout the column)
$ws.Cells.Item(73, 4).NumberFormat = $ws.Cells.Item(72, 4).NumberFormat
for ($col = 1; $col -le 20; $col++) {
    if ($col -eq 4 -or $col -eq 12 -or $col -eq 13 -or $col -eq 14 -or $col -eq 15 -or $col -eq 16 -or $col -eq 19) {
        continue
    }
    $ws.Cells.Item(73, $col).Value2 = $ws.Cells.Item(72, $col).Value2
}

# --- Step 3: row 44 gets a brand-new latest observation; N,O,P,S stay as before ---
$ws.Cells.Item(44, 4).Value  = 44438   # D: Fecha
$ws.Cells.Item(44, 13).Value = 130     # M: Calibre/tamaño

